# Daily update: prepend the newest day's record (2021-12-17 / serial 44547)
# for "Macroferia Regional de Talca - Arveja Verde" at row 10, shifting all
# existing data rows (10-84) down by one (to 11-85).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 10; this pushes former rows 10..84 to 11..85
# and carries the row's cell formatting (date style on column D) down with it.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new day's data.
$ws.Cells.Item(10, 1).Value  = 5
$ws.Cells.Item(10, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(10, 3).Value  = "Maule"
$ws.Cells.Item(10, 4).Value  = 44547
$ws.Cells.Item(10, 5).Value  = 7
$ws.Cells.Item(10, 6).Value  = 100112022
$ws.Cells.Item(10, 7).Value  = "Arveja Verde"
$ws.Cells.Item(10, 8).Value  = "Sin especificar"
$ws.Cells.Item(10, 9).Value  = "Primera"
$ws.Cells.Item(10, 10).Value = 300
$ws.Cells.Item(10, 11).Value = 15000
$ws.Cells.Item(10, 12).Value = 15000
$ws.Cells.Item(10, 13).Value = 15000
$ws.Cells.Item(10, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(10, 15).Value = "Carahue"
$ws.Cells.Item(10, 16).Value = 600
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
